# "Kleine foutjes in de dataset verbeterd."
# Fix small mistakes in the dataset:
# 1. Row 8 (Knop) last name should be "Knops"
# 2. Row 17 (Slob) Stropdas (necktie) value should be 1, not 0
# 3. Update the active selection to A19 (as left by the editor)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Knops"
$ws.Range("C17").Value = 1

$ws.Range("A19").Select()
